{"js": "// The document contains a single table of simple arithmetic expressions\n// (e.g. \"56-0=\", \"76+20=\"). The edit replaces the text of every cell with\n// a new expression while leaving every other part of the document\n// (formatting, paragraph/run properties, table structure, the date\n// heading, etc.) untouched.\n//\n// We read the table's current grid of values and overwrite it in one\n// shot with the new grid (same row/column shape), which updates each\n// cell's existing run text in place without touching formatting.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New text for every cell, row-major (20 rows x 5 columns), in the exact\n// order the cells appear in the document.\nconst newValues = [\n  [\"23+33=\", \"33+21=\", \"16+45=\", \"44-16=\", \"63-51=\"],\n  [\"43-19=\", \"74-24=\", \"88-74=\", \"36-35=\", \"46+32=\"],\n  [\"93-81=\", \"85+3=\", \"3+36=\", \"20+42=\", \"1+28=\"],\n  [\"60+8=\", \"61-54=\", \"9+21=\", \"88-9=\", \"85-0=\"],\n  [\"68-64=\", \"21+24=\", \"8+72=\", \"80-33=\", \"20+45=\"],\n  [\"66-8=\", \"59-34=\", \"70-59=\", \"99-52=\", \"37+2=\"],\n  [\"84-37=\", \"48-28=\", \"15+84=\", \"37+52=\", \"58-22=\"],\n  [\"44-10=\", \"63-23=\", \"70-5=\", \"98-63=\", \"11+51=\"],\n  [\"74-42=\", \"39-8=\", \"18+39=\", \"18+44=\", \"33+8=\"],\n  [\"33+45=\", \"2+17=\", \"77+2=\", \"87-15=\", \"85-71=\"],\n  [\"66-41=\", \"65-65=\", \"14+6=\", \"79+18=\", \"28-27=\"],\n  [\"3+82=\", \"93-75=\", \"49-22=\", \"94-39=\", \"52+9=\"],\n  [\"53-13=\", \"77+10=\", \"88-13=\", \"0+36=\", \"57-5=\"],\n  [\"88-25=\", \"18+47=\", \"66-18=\", \"10+88=\", \"78-7=\"],\n  [\"81-68=\", \"70-23=\", \"68-68=\", \"25+31=\", \"41+21=\"],\n  [\"1+75=\", \"66-52=\", \"58-2=\", \"18+7=\", \"75-26=\"],\n  [\"40+53=\", \"14+43=\", \"95-1=\", \"11-0=\", \"67+16=\"],\n  [\"83-21=\", \"76-33=\", \"2+16=\", \"83-0=\", \"2+73=\"],\n  [\"53-37=\", \"5+85=\", \"38-5=\", \"78-57=\", \"92-85=\"],\n  [\"8-0=\", \"95-12=\", \"82+8=\", \"89-50=\", \"5+42=\"],\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# The document contains a single table of simple arithmetic expressions\n# (e.g. \"56-0=\", \"76+20=\"). The edit replaces the text of every cell with\n# a new expression while leaving every other part of the document\n# (formatting, paragraph/run properties, table structure, the date\n# heading, etc.) untouched.\n#\n# Walk the table cell by cell (row-major, matching document order) and\n# overwrite each cell's Range.Text with the new expression. Setting\n# Range.Text on a table cell replaces only the cell's content, leaving\n# the end-of-cell marker (and the run/paragraph formatting already on\n# that text) intact.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# New text for every cell, row-major (20 rows x 5 columns), in the exact\n# order the cells appear in the document.\n$newValues = @(\n    @(\"23+33=\",\"33+21=\",\"16+45=\",\"44-16=\",\"63-51=\"),\n    @(\"43-19=\",\"74-24=\",\"88-74=\",\"36-35=\",\"46+32=\"),\n    @(\"93-81=\",\"85+3=\",\"3+36=\",\"20+42=\",\"1+28=\"),\n    @(\"60+8=\",\"61-54=\",\"9+21=\",\"88-9=\",\"85-0=\"),\n    @(\"68-64=\",\"21+24=\",\"8+72=\",\"80-33=\",\"20+45=\"),\n    @(\"66-8=\",\"59-34=\",\"70-59=\",\"99-52=\",\"37+2=\"),\n    @(\"84-37=\",\"48-28=\",\"15+84=\",\"37+52=\",\"58-22=\"),\n    @(\"44-10=\",\"63-23=\",\"70-5=\",\"98-63=\",\"11+51=\"),\n    @(\"74-42=\",\"39-8=\",\"18+39=\",\"18+44=\",\"33+8=\"),\n    @(\"33+45=\",\"2+17=\",\"77+2=\",\"87-15=\",\"85-71=\"),\n    @(\"66-41=\",\"65-65=\",\"14+6=\",\"79+18=\",\"28-27=\"),\n    @(\"3+82=\",\"93-75=\",\"49-22=\",\"94-39=\",\"52+9=\"),\n    @(\"53-13=\",\"77+10=\",\"88-13=\",\"0+36=\",\"57-5=\"),\n    @(\"88-25=\",\"18+47=\",\"66-18=\",\"10+88=\",\"78-7=\"),\n    @(\"81-68=\",\"70-23=\",\"68-68=\",\"25+31=\",\"41+21=\"),\n    @(\"1+75=\",\"66-52=\",\"58-2=\",\"18+7=\",\"75-26=\"),\n    @(\"40+53=\",\"14+43=\",\"95-1=\",\"11-0=\",\"67+16=\"),\n    @(\"83-21=\",\"76-33=\",\"2+16=\",\"83-0=\",\"2+73=\"),\n    @(\"53-37=\",\"5+85=\",\"38-5=\",\"78-57=\",\"92-85=\"),\n    @(\"8-0=\",\"95-12=\",\"82+8=\",\"89-50=\",\"5+42=\")\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$r-1][$c-1]\n    }\n}\n"}
